# regen sval data to filter save games
# Update the B:G numeric columns (rows 2-12) of Sheet1 with the regenerated
# statistic values. Column A (dates) and column F (win flag) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    3  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4  = @(3.272327238179451, 1.626987699542094, 3.223369029078222,  0.5333859586016987, 8.656069925401464)
    5  = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 16.86649396021207)
    6  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    7  = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 2.213936997104367)
    8  = @(1.445647641019636, 0.3048912486333797, 18.71679738969934, 13.86384647080068, 34.33118275015303)
    9  = @(1.445647641019636, 1.626987699542094, 18.71679738969934,  0.5333859586016987, 22.32281868886277)
    10 = @(1.445647641019636, 1.626987699542094, 3.223369029078222,  0.5333859586016987, 6.82939032824165)
    11 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    12 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
